$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full replacement table (17 rows x 6 cols).
# New layout inserts a "groupID" column before "groupName" and adds a
# "shortcutBg" column after "link"; link values become full https URLs
# and each shortcut gets its own accent color.
$data = New-Object 'object[,]' 17,6

# Header row
$data[0,0] = 'groupID'
$data[0,1] = 'groupName'
$data[0,2] = 'groupBg'
$data[0,3] = 'shortcutName'
$data[0,4] = 'link'
$data[0,5] = 'shortcutBg'

# group-1 / MAIN rows
$data[1,0]  = 'group-1'; $data[1,1]  = 'MAIN'; $data[1,2]  = '#6c8d91'; $data[1,3]  = 'Google';                $data[1,4]  = 'https://www.google.com/';    $data[1,5]  = '#F39C12'
$data[2,0]  = 'group-1'; $data[2,1]  = 'MAIN'; $data[2,2]  = '#6c8d91'; $data[2,3]  = 'YouTube';               $data[2,4]  = 'https://www.youtube.com/';   $data[2,5]  = '#E74C3C'
$data[3,0]  = 'group-1'; $data[3,1]  = 'MAIN'; $data[3,2]  = '#6c8d91'; $data[3,3]  = 'Facebook';              $data[3,4]  = 'https://www.facebook.com/';  $data[3,5]  = '#8E44AD'
$data[4,0]  = 'group-1'; $data[4,1]  = 'MAIN'; $data[4,2]  = '#6c8d91'; $data[4,3]  = 'Instagram';             $data[4,4]  = 'https://www.instagram.com/'; $data[4,5]  = '#3498DB'
$data[5,0]  = 'group-1'; $data[5,1]  = 'MAIN'; $data[5,2]  = '#6c8d91'; $data[5,3]  = 'WhatsApp';              $data[5,4]  = 'https://www.whatsapp.com/';  $data[5,5]  = '#2ECC71'
$data[6,0]  = 'group-1'; $data[6,1]  = 'MAIN'; $data[6,2]  = '#6c8d91'; $data[6,3]  = 'X (formerly Twitter)';  $data[6,4]  = 'https://www.x.com/';         $data[6,5]  = '#D35400'
$data[7,0]  = 'group-1'; $data[7,1]  = 'MAIN'; $data[7,2]  = '#6c8d91'; $data[7,3]  = 'Wikipedia';             $data[7,4]  = 'https://www.wikipedia.org/'; $data[7,5]  = '#C0392B'
$data[8,0]  = 'group-1'; $data[8,1]  = 'MAIN'; $data[8,2]  = '#6c8d91'; $data[8,3]  = 'ChatGPT';               $data[8,4]  = 'https://www.chatgpt.com/';   $data[8,5]  = '#27AE60'
$data[9,0]  = 'group-1'; $data[9,1]  = 'MAIN'; $data[9,2]  = '#6c8d91'; $data[9,3]  = 'Reddit';                $data[9,4]  = 'https://www.reddit.com/';    $data[9,5]  = '#2980B9'
$data[10,0] = 'group-1'; $data[10,1] = 'MAIN'; $data[10,2] = '#6c8d91'; $data[10,3] = 'Yahoo';                 $data[10,4] = 'https://www.yahoo.com/';     $data[10,5] = '#8E44AD'

# group-2 / SUPP rows (no groupBg value for this group)
$data[11,0] = 'group-2'; $data[11,1] = 'SUPP'; $data[11,2] = $null; $data[11,3] = 'Amazon';    $data[11,4] = 'https://www.amazon.com/';    $data[11,5] = '#F1C40F'
$data[12,0] = 'group-2'; $data[12,1] = 'SUPP'; $data[12,2] = $null; $data[12,3] = 'LinkedIn';  $data[12,4] = 'https://www.linkedin.com/';  $data[12,5] = '#E67E22'
$data[13,0] = 'group-2'; $data[13,1] = 'SUPP'; $data[13,2] = $null; $data[13,3] = 'Netflix';   $data[13,4] = 'https://www.netflix.com/';   $data[13,5] = '#1ABC9C'
$data[14,0] = 'group-2'; $data[14,1] = 'SUPP'; $data[14,2] = $null; $data[14,3] = 'eBay';      $data[14,4] = 'https://www.ebay.com/';      $data[14,5] = '#34495E'
$data[15,0] = 'group-2'; $data[15,1] = 'SUPP'; $data[15,2] = $null; $data[15,3] = $null;       $data[15,4] = $null;                        $data[15,5] = '#7D3C98'
$data[16,0] = 'group-2'; $data[16,1] = 'SUPP'; $data[16,2] = $null; $data[16,3] = 'Pinterest'; $data[16,4] = 'https://www.pinterest.com/'; $data[16,5] = '#F39C12'

# Clear the old used range first so no stray formatting/values linger,
# then write the new 17x6 table in one shot.
$ws.Range("A1:F17").Clear()
$ws.Range("A1:F17").Value = $data

# Widen the new "link" column (E) to fit the full URLs.
$ws.Columns.Item(5).ColumnWidth = 37.022135416666664

# Match the saved selection/active cell.
$ws.Range("J11").Select() | Out-Null
